$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the whole B2:K21 block with the background/filler score first.
$ws.Range("B2:K21").Value = -16.62100729546654

# Then overwrite the specific residue-pair cells with their real PSSM scores.
$ws.Range("C2").Value = 1.987626548075904
$ws.Range("C4").Value = 2.157599271101367
$ws.Range("D4").Value = 2.285768823947909
$ws.Range("F4").Value = 3.398975640520367
$ws.Range("H4").Value = 1.584523479579194
$ws.Range("J4").Value = 2.341846358090495
$ws.Range("C5").Value = 1.905419348607962
$ws.Range("G5").Value = 2.872920452419557
$ws.Range("B7").Value = 2.580506442036932
$ws.Range("E8").Value = 1.763982977501352
$ws.Range("B9").Value = 3.809242077042278
$ws.Range("I10").Value = 4.321914496780905
$ws.Range("K10").Value = 1.630738606322347
$ws.Range("E11").Value = 2.916076468085719
$ws.Range("G11").Value = 2.609691027351073
$ws.Range("K11").Value = 1.747025908734648
$ws.Range("E13").Value = 2.373735603220104
$ws.Range("J13").Value = 2.231391149205259
$ws.Range("K13").Value = 1.586335964702213
$ws.Range("D14").Value = 1.313564721815956
$ws.Range("K14").Value = 2.249625327631868
$ws.Range("D15").Value = 1.160614958517817
$ws.Range("J16").Value = 2.261614859309395
$ws.Range("C17").Value = 1.895065886697096
$ws.Range("D17").Value = 2.184046925548622
$ws.Range("H17").Value = 1.057242084825161
$ws.Range("J17").Value = 1.437433719274042
$ws.Range("H18").Value = 0.9367929279119046
$ws.Range("J18").Value = 1.449328569398794
$ws.Range("D19").Value = 1.56945883748148
$ws.Range("H19").Value = 1.949667456309237
$ws.Range("C20").Value = 0.9128166966820443
$ws.Range("D20").Value = 1.531242881491655
$ws.Range("F20").Value = 3.240505242152876
$ws.Range("H20").Value = 2.241899404267347
$ws.Range("K20").Value = 2.533071579172516
$ws.Range("C21").Value = 1.153168475843
$ws.Range("E21").Value = 1.953368185149882
$ws.Range("G21").Value = 2.716033080590695
$ws.Range("H21").Value = 2.141598589663956
